$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two re-ordered rows)

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '61.755.95'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('E2').Style = "Normal"

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.405.54'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E3').Style = "Normal"

# Row 4
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('E4').Style = "Normal"

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '411.32'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('E5').Style = "Normal"

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '130.51'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('E6').Style = "Normal"

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.618'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -2.21%  '
$ws.Range('E7').Style = "Normal"

# Row 8
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E8').Style = "Normal"

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.723'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('E9').Style = "Normal"

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.133'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -6.22%  '
$ws.Range('E10').Style = "Normal"

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.66'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('E11').Style = "Normal"

# Row 12
$ws.Range('B12').NumberFormat = "@"
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('B12').Style = "Normal"
$ws.Range('C12').NumberFormat = "@"
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C12').Style = "Normal"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '3.951.95'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E12').Style = "Normal"

# Row 13
$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('B13').Style = "Normal"
$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C13').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.09'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('E13').Style = "Normal"

# Row 14
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('E14').Style = "Normal"

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000207'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('E15').Style = "Normal"

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.38'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('E16').Style = "Normal"

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.404.92'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E17').Style = "Normal"

# Row 18
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.61%  '
$ws.Range('E18').Style = "Normal"

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.25'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('E19').Style = "Normal"

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '61.805.76'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('E20').Style = "Normal"

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '477.04'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +18.46%  '
$ws.Range('E21').Style = "Normal"

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '89.70'
$ws.Range('D22').Style = "Normal"

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.24'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.39%  '
$ws.Range('E23').Style = "Normal"

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.13'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('E24').Style = "Normal"

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.29'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.27%  '
$ws.Range('E25').Style = "Normal"

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.69'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +12.79%  '
$ws.Range('E26').Style = "Normal"

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '32.92'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('E27').Style = "Normal"

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '4.76'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E28').Style = "Normal"

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +5.12%  '
$ws.Range('E29').Style = "Normal"

# Row 30
$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('B30').Style = "Normal"
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C30').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '11.86'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('E30').Style = "Normal"

# Row 31
$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('B31').Style = "Normal"
$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C31').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.65'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('E31').Style = "Normal"

# Row 32
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -1.97%  '
$ws.Range('E32').Style = "Normal"

# Row 33
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -5.00%  '
$ws.Range('E33').Style = "Normal"

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '40.92'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('E34').Style = "Normal"

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('E35').Style = "Normal"

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '56.80'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +5.26%  '
$ws.Range('E36').Style = "Normal"

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0485'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.61%  '
$ws.Range('E37').Style = "Normal"

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.27%  '
$ws.Range('E38').Style = "Normal"

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.04'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +4.74%  '
$ws.Range('E39').Style = "Normal"

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.328'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +5.91%  '
$ws.Range('E40').Style = "Normal"

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '147.86'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +5.26%  '
$ws.Range('E41').Style = "Normal"

# Row 42
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.37%  '
$ws.Range('E42').Style = "Normal"

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.33'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('E43').Style = "Normal"

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.06'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +5.18%  '
$ws.Range('E44').Style = "Normal"

# Row 45
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +8.31%  '
$ws.Range('E45').Style = "Normal"

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.23'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +4.83%  '
$ws.Range('E46').Style = "Normal"

# Row 47
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('B47').Style = "Normal"
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('C47').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.33'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +19.36%  '
$ws.Range('E47').Style = "Normal"

# Row 48
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'Celestia'
$ws.Range('B48').Style = "Normal"
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('C48').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '16.49'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('E48').Style = "Normal"

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '22.02'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.33%  '
$ws.Range('E49').Style = "Normal"

# Row 50
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'Cronos'
$ws.Range('B50').Style = "Normal"
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C50').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.141'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +7.65%  '
$ws.Range('E50').Style = "Normal"

# Row 51
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('B51').Style = "Normal"
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('C51').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '112.05'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +14.10%  '
$ws.Range('E51').Style = "Normal"
